# "Update pure guess model" - add the EXP(-0.5*C) guess-model formulas
# back into column D for rows 9:13 (Pure Guess / Pure Intrusion /
# Intrusion + Guess / Temporal / Spatiotemporal), rolling back to the
# older calculation that had been blanked out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Formula  = "=EXP(-0.5*C9)"
$ws.Range("D10").Formula = "=EXP(-0.5*C10)"
$ws.Range("D11").Formula = "=EXP(-0.5*C11)"
$ws.Range("D12").Formula = "=EXP(-0.5*C12)"
$ws.Range("D13").Formula = "=EXP(-0.5*C13)"

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("H28").Select()
